$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Book of Gems Megaways Free for
#    an Ancient Egypt Adventure").
# ----------------------------------------------------------------------
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titleIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Play Book of Gems Megaways Free for an Ancient Egypt Adventure*") {
        $titleIndex = $i
        break
    }
}
$titlePara = $d.Paragraphs($titleIndex)
[void]$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs($titleIndex + 1)
$metaPara.Range.Style = "Normal"
$metaXml = "<w:p $wNs>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
    "<w:r><w:t>: Read our review of Book of Gems Megaways slot play free. Enjoy high RTP, free spins with expanding symbol, and impressive graphics set in Ancient Egypt.</w:t></w:r>" +
    "</w:p>"
[void]$metaPara.Range.InsertXML($metaXml)

# ----------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Book of Gems Megaways Free for an
#    Ancient Egypt Adventure" paragraph that was near the end of the
#    document, right before the final (italic) paragraph.
# ----------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Play Book of Gems Megaways Free for an Ancient Egypt Adventure`r") {
        [void]$p.Range.Delete()
        break
    }
}

# ----------------------------------------------------------------------
# 3) Replace the text of the paragraph that used to hold the meta
#    description ("Read our review of ...") with the image-generation
#    prompt, keeping the italic formatting.
# ----------------------------------------------------------------------
$promptXml = "<w:p $wNs>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:i/></w:rPr><w:t>Prompt: Create a feature image for Book of Gems Megaways that portrays a happy Maya warrior with glasses in a cartoon style. The image should showcase the warrior holding a golden book and standing in front of a desert landscape with pyramids in the background. The warrior should be dressed in traditional Mayan clothing with an assortment of precious stones and gems adorning their attire. The cartoon style of the image should be vibrant and colorful, with the warrior looking excited and thrilled to be playing the game. The background should be a warm and inviting sunset color, with a sandstorm brewing in the distance to create an ominous yet exciting atmosphere. Overall, the feature image should showcase the game's theme of adventure and excitement while highlighting the warrior's enthusiasm and love for playing online slot games.</w:t></w:r>" +
    "</w:p>"

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Read our review of Book of Gems Megaways*") {
        [void]$p.Range.InsertXML($promptXml)
        break
    }
}
